$d = $word.ActiveDocument
$d.TrackRevisions = $false
$d.Content.LanguageID = 1033
